# Apply updated "dSF" (column F) values for the 2023 castillo_luis save-data sheet.
# These values represent a recalculated/repulled metric (see commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = -4
    3  = -5
    4  = 4
    5  = 1
    6  = 1
    7  = 3
    8  = 3
    9  = 8
    10 = 0
    11 = 7
    12 = 0
    13 = 3
    14 = -1
    15 = -3
    16 = -3
    17 = 8
    18 = -1
    19 = 3
    20 = -3
    21 = -2
    22 = -2
    23 = -1
    24 = 4
    25 = 6
    26 = -3
    27 = -2
    28 = -1
    30 = -2
    31 = 1
    32 = -1
    33 = 8
    35 = 4
    36 = -3
    38 = -2
    39 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
